# Auto-generated Excel COM-interop script applying the Masamune_Profits diff
# Updates currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns (H-N)
# across several sheets, matching the scheduled-runner data refresh described
# in the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 141.77777
$ws.Range("I33").Value = 132.5
$ws.Range("K33").Value = 132.5
$ws.Range("M33").Value = 96.5

$ws.Range("H112").Value = 1135.1428
$ws.Range("J112").Value = 1135.1428
$ws.Range("L112").Value = 3405.4284
$ws.Range("N112").Value = -5621.428400000001

$ws.Range("H116").Value = 5964.5654
$ws.Range("I116").Value = 3033.7
$ws.Range("J116").Value = 8219.076999999999
$ws.Range("K116").Value = 3033.7
$ws.Range("L116").Value = 8219.076999999999
$ws.Range("M116").Value = 408.3000000000002
$ws.Range("N116").Value = -15103.077

$ws.Range("H123").Value = 39971.332
$ws.Range("J123").Value = 39971.332
$ws.Range("L123").Value = 39971.332
$ws.Range("N123").Value = -49771.332

$ws.Range("H132").Value = 15592.667
$ws.Range("I132").Value = 2518.5254
$ws.Range("J132").Value = 125789
$ws.Range("K132").Value = 7555.5762
$ws.Range("L132").Value = 377367
$ws.Range("M132").Value = -5025.5762
$ws.Range("N132").Value = -382427

$ws.Range("H135").Value = 17858304
$ws.Range("I135").Value = 1231.3636
$ws.Range("J135").Value = 83334240
$ws.Range("K135").Value = 11082.2724
$ws.Range("L135").Value = 750008160
$ws.Range("M135").Value = -8547.2724
$ws.Range("N135").Value = -750013230

$ws.Range("H137").Value = 3073.9487
$ws.Range("I137").Value = 1235.3572
$ws.Range("J137").Value = 3476.1406
$ws.Range("K137").Value = 3706.0716
$ws.Range("L137").Value = 10428.4218
$ws.Range("M137").Value = -1156.0716
$ws.Range("N137").Value = -15528.4218

$ws.Range("H138").Value = 2573.4524
$ws.Range("I138").Value = 2023.2285
$ws.Range("J138").Value = 2966.4695
$ws.Range("K138").Value = 6069.6855
$ws.Range("L138").Value = 8899.408500000001
$ws.Range("M138").Value = -929.6854999999996
$ws.Range("N138").Value = -19179.4085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 24957.5
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 24957.5
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 24957.5
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -25891.5

$ws.Range("H61").Value = 2720.25
$ws.Range("I61").Value = 2662.72
$ws.Range("J61").Value = 3199.6667
$ws.Range("K61").Value = 2662.72
$ws.Range("L61").Value = 3199.6667
$ws.Range("M61").Value = -2450.72
$ws.Range("N61").Value = -3623.6667

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H110").Value = 1739.4584
$ws.Range("I110").Value = 1715.7727
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1715.7727
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 329.2273
$ws.Range("N110").Value = -6090

$ws.Range("H136").Value = 2720.25
$ws.Range("I136").Value = 2662.72
$ws.Range("J136").Value = 3199.6667
$ws.Range("K136").Value = 7988.16
$ws.Range("L136").Value = 9599.000100000001
$ws.Range("M136").Value = -5438.16
$ws.Range("N136").Value = -14699.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2155.04
$ws.Range("I105").Value = 1947.0834
$ws.Range("J105").Value = 2347
$ws.Range("K105").Value = 1947.0834
$ws.Range("L105").Value = 2347
$ws.Range("M105").Value = -200.0834
$ws.Range("N105").Value = -5841

$ws.Range("H107").Value = 2998.5715
$ws.Range("I107").Value = 2998
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2998
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -1078
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 4742.4
$ws.Range("I134").Value = 4521
$ws.Range("K134").Value = 13563
$ws.Range("M134").Value = -11028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 461.85715
$ws.Range("J22").Value = 102
$ws.Range("L22").Value = 102
$ws.Range("N22").Value = -802

$ws.Range("H31").Value = 5728.709
$ws.Range("I31").Value = 2706.9583
$ws.Range("J31").Value = 8068.129
$ws.Range("K31").Value = 2706.9583
$ws.Range("L31").Value = 8068.129
$ws.Range("M31").Value = -2411.9583
$ws.Range("N31").Value = -8658.129000000001

$ws.Range("H34").Value = 5728.709
$ws.Range("I34").Value = 2706.9583
$ws.Range("J34").Value = 8068.129
$ws.Range("K34").Value = 2706.9583
$ws.Range("L34").Value = 8068.129
$ws.Range("M34").Value = -2504.9583
$ws.Range("N34").Value = -8472.129000000001

$ws.Range("H132").Value = 49140.266
$ws.Range("I132").Value = 1880.2667
$ws.Range("J132").Value = 96400.266
$ws.Range("K132").Value = 5640.800099999999
$ws.Range("L132").Value = 289200.798
$ws.Range("M132").Value = -3110.800099999999
$ws.Range("N132").Value = -294260.798

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1739.75
$ws.Range("J117").Value = 2816.6667
$ws.Range("L117").Value = 8450.000100000001
$ws.Range("N117").Value = -15334.0001

$ws.Range("H131").Value = 2306.1829
$ws.Range("I131").Value = 20625.2
$ws.Range("J131").Value = 1116.6364
$ws.Range("K131").Value = 61875.60000000001
$ws.Range("L131").Value = 3349.9092
$ws.Range("M131").Value = -56835.60000000001
$ws.Range("N131").Value = -13429.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6000

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H105").Value = 44397
$ws.Range("J105").Value = 44397
$ws.Range("L105").Value = 44397
$ws.Range("N105").Value = -51385

$ws.Range("H122").Value = 904
$ws.Range("I122").Value = 751
$ws.Range("J122").Value = 1210
$ws.Range("K122").Value = 2253
$ws.Range("L122").Value = 3630
$ws.Range("M122").Value = 197
$ws.Range("N122").Value = -8530

$ws.Range("H126").Value = 6853.8096
$ws.Range("I126").Value = 12001.2
$ws.Range("J126").Value = 2174.3635
$ws.Range("K126").Value = 36003.60000000001
$ws.Range("L126").Value = 6523.0905
$ws.Range("M126").Value = -33533.60000000001
$ws.Range("N126").Value = -11463.0905

$ws.Range("H132").Value = 2484.1282
$ws.Range("I132").Value = 1729.2222
$ws.Range("J132").Value = 4182.6665
$ws.Range("K132").Value = 5187.6666
$ws.Range("L132").Value = 12547.9995
$ws.Range("M132").Value = -2657.6666
$ws.Range("N132").Value = -17607.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 31450.455
$ws.Range("J123").Value = 31450.455
$ws.Range("L123").Value = 31450.455
$ws.Range("N123").Value = -41250.455

$ws.Range("H132").Value = 4443.4053
$ws.Range("I132").Value = 5484.5386
$ws.Range("J132").Value = 3879.4583
$ws.Range("K132").Value = 16453.6158
$ws.Range("L132").Value = 11638.3749
$ws.Range("M132").Value = -13923.6158
$ws.Range("N132").Value = -16698.3749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1607.2727
$ws.Range("I81").Value = 1568
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3136
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -2075
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 1607.2727
$ws.Range("I84").Value = 1568
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 15680
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -10376
$ws.Range("N84").Value = -30608

$ws.Range("H132").Value = 2904.6
$ws.Range("I132").Value = 2089.4
$ws.Range("J132").Value = 3719.8
$ws.Range("K132").Value = 6268.200000000001
$ws.Range("L132").Value = 11159.4
$ws.Range("M132").Value = -3738.200000000001
$ws.Range("N132").Value = -16219.4

$ws.Range("H136").Value = 19975.908
$ws.Range("I136").Value = 39653.652
$ws.Range("K136").Value = 118960.956
$ws.Range("M136").Value = -116410.956
